# Add the three new character styles used by the edited paragraphs.
$d = $word.ActiveDocument

$styleDate = $d.Styles.Add("GaNStyle", 2)
$styleDate.Font.Name = "Calibri"
$styleDate.Font.Size = 14

$styleParagraph = $d.Styles.Add("GaNParagraph", 2)
$styleParagraph.Font.Name = "Calibri"
$styleParagraph.Font.Size = 10

$styleLinks = $d.Styles.Add("GaNLinks", 2)
$styleLinks.Font.Name = "Calibri"
$styleLinks.Font.Bold = $true
$styleLinks.Font.Color = 8388608
$styleLinks.Font.Size = 9.5
$styleLinks.Font.Underline = 1

# The four "observation window" date-range paragraphs: apply GaNStyle and
# append the missing trailing period to the sentence.
$oldDates = "Herkuleen tähtikuvio havainnointijaksot vuonna 2022: 13.-22.6., 12.-21.7., 10.-19.8"
$newDates = "Herkuleen tähtikuvio havainnointijaksot vuonna 2022: 13.-22.6., 12.-21.7., 10.-19.8."

$rng = $d.Content
$guard = 0
while ($rng.Find.Execute($oldDates, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Text = $newDates
    $rng.Collapse(0)
    $rng = $d.Range($rng.End, $d.Content.End)
    $guard = $guard + 1
    if ($guard -ge 20) { break }
}

# The introductory "Osallistut maailmanlaajuiseen..." paragraph: apply GaNParagraph.
$introText = "Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpiä näkyvissä olevia tähtiä keinona mitata valonsaastetta tietyssä paikassa. Paikallistamalla ja tarkkailemalla Herkuleen tähtikuvio miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot päivittyvät heti verkossa olevaan tietokantaan, ja näin saadaan käsitys siitä minkä verran taivaan tähdistä on missäkin nähtävissä."

$rng2 = $d.Content
if ($rng2.Find.Execute($introText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng2.Style = "GaNParagraph"
}

# The "Tämän oppaan kartat piirsi..." credit line: apply GaNLinks.
$creditText = "Tämän oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng3 = $d.Content
if ($rng3.Find.Execute($creditText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng3.Style = "GaNLinks"
}

Write-Output "Done."
